# Update "02 Software-Projekt Tätigkeitsdokumentation.xlsx"
# Adds a new log entry (row 33) to the Tätigkeitsdokumentation sheet and
# moves the active selection/viewport further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tätigkeitsdokumentation")

# New row of data: Commit | Datei | Abschnitt/Klasse | Quelltextzeilen | Schwierigkeitsgrad | Docs | Programmierer/-in
$ws.Range("A33").Value = "5 commits"
$ws.Range("B33").Value = "multiple files"
$ws.Range("C33").Value = "smaller updates"
$ws.Range("D33").Value = 56
$ws.Range("E33").Value = 1
$ws.Range("F33").Value = "inline 100%"
$ws.Range("G33").Value = "Giesbrt"

# Update the view: scroll so row 13 is at the top and select B34.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B34").Select()
